# Atualização de bases das ligas, do dia: 13-06-2024 às 19:35
# Swap the B:AD content (everything except the id column A) between the
# pairs of rows that had their data mixed up: (128,129), (143,145), (148,149)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($ws, $row1, $row2) {
    $addr1 = "B" + $row1 + ":AD" + $row1
    $addr2 = "B" + $row2 + ":AD" + $row2

    $range1 = $ws.Range($addr1)
    $range2 = $ws.Range($addr2)

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}

Swap-RowData $ws 128 129
Swap-RowData $ws 143 145
Swap-RowData $ws 148 149
